# Applies the "added citation, moved down sources" edit to slide 1 of the
# Wikipedia word2vec poster.
#
#   1. Moves the "Sources" section header textbox down.
#   2. Repurposes the (previously unused) "John A. Paulson..." textbox into a
#      small citation line ("code.google.com/p/word2vec/source/") set in
#      Courier New, repositioned/resized next to the Sources header, with
#      word-wrap turned on and center alignment cleared.
#   3. Moves the two source-citation pictures further down to make room.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "Sources" section header: shift down ---------------------------
$sourcesHeader = $s.Shapes.Item(10)
$sourcesHeader.Top = 33676050 / 12700

# --- 2. Citation textbox (formerly "John A. Paulson..." textbox) -------
$citation = $s.Shapes.Item(16)

$citation.Left = 19404013 / 12700
$citation.Top = 35955165 / 12700
$citation.Width = 5101907 / 12700
$citation.Height = 300846 / 12700

$citation.TextFrame.WordWrap = -1

$tf = $citation.TextFrame
$tr = $tf.TextRange
$tr.ParagraphFormat.Alignment = 1

$tr.Text = "code.google.com"
$tr.Font.Name = "Courier New"
$tr.Font.NameComplexScript = "Courier New"
$tr.Font.Size = 14

$run2 = $tf.TextRange.InsertAfter("/p/word2vec/source/")
$run2.Font.Name = "Courier New"
$run2.Font.NameComplexScript = "Courier New"
$run2.Font.Size = 14

# --- 3. Move the two source-citation pictures further down -------------
$pic21 = $s.Shapes.Item(27)
$pic21.Top = 36434332 / 12700

$pic23 = $s.Shapes.Item(28)
$pic23.Top = 37162996 / 12700
